$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 12, pushing the existing rows 12-23 down to 13-24
# (dimension grows from A1:R23 to A1:R24).
$ws.Rows.Item(12).EntireRow.Insert()

# Populate the newly inserted row 12 with the new weekly price record.
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C12").Value = "Los Lagos"
$ws.Range("D12").Value = 44789
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 100112013
$ws.Range("G12").Value = "Alcachofa"
$ws.Range("H12").Value = "Madrigal"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 16000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 16000
$ws.Range("N12").Value = "$/caja 40 unidades"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 400
$ws.Range("Q12").Value = 40
$ws.Range("R12").Value = "Hortaliza"
